# Adapt column header formatting to respective input file names (#7):
#   "<header>_old" -> "<header>_FV2410"
#   "<header>_new" -> "<header>_FV2504"
# Then wrap the data range in a structured Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1, columns A:U) ------------------------
# Columns A-J (1-10) carry the "_old" suffix  -> replace with "_FV2410"
# Column  K    (11)  is "diff"                -> left untouched
# Columns L-U (12-21) carry the "_new" suffix -> replace with "_FV2504"
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Value()
    if ($header -like "*_old") {
        $base = $header.Substring(0, $header.Length - 4)
        $cell.Value = "$($base)_FV2410"
    } elseif ($header -like "*_new") {
        $base = $header.Substring(0, $header.Length - 4)
        $cell.Value = "$($base)_FV2504"
    }
}

# --- 2. Turn the data range into a structured Table ---------------------
$dataRange = $ws.Range("A1:U73")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# --- 3. Freeze the header row (row 1) -----------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
